# Add Thingiverse links to the BOM "Link" column (K) for rows 26-38.
# Commit message: "Link to Thingiverse Added to BOM"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 26-33 and 38 already carry style s="30" (matches the rest of the
# table); rows 35-37 currently use style s="17" and must be switched to the
# same s="30" look used elsewhere in the table. Row 34 is the section-header
# row: its K cell must match the rest of that header row (style s="45",
# same as J34) rather than the special corner style s="46" it still has.

$ws.Range("J34").Copy()
$ws.Range("K34").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("K26").Copy()
$ws.Range("K35").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("K26").Copy()
$ws.Range("K36").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("K26").Copy()
$ws.Range("K37").PasteSpecial(-4122)   # xlPasteFormats

# Now fill in the Thingiverse link text for every BOM row in this block.
$ws.Range("K26").Value = "https://www.thingiverse.com/thing:5173650"
$ws.Range("K27").Value = "https://www.thingiverse.com/thing:5173651"
$ws.Range("K28").Value = "https://www.thingiverse.com/thing:5173652"
$ws.Range("K29").Value = "https://www.thingiverse.com/thing:5173653"
$ws.Range("K30").Value = "https://www.thingiverse.com/thing:5173654"
$ws.Range("K31").Value = "https://www.thingiverse.com/thing:5173655"
$ws.Range("K32").Value = "https://www.thingiverse.com/thing:5173656"
$ws.Range("K33").Value = "https://www.thingiverse.com/thing:5173657"
$ws.Range("K35").Value = "https://www.thingiverse.com/thing:5173659"
$ws.Range("K36").Value = "https://www.thingiverse.com/thing:5173660"
$ws.Range("K37").Value = "https://www.thingiverse.com/thing:5173661"
$ws.Range("K38").Value = "https://www.thingiverse.com/thing:5173662"

# Leave the selection on K35, matching where the author ended up editing -
# this also drops the stale topLeftCell="A8" scroll position from the
# saved view since K35 is already in the default viewport.
$ws.Range("K35").Select()
